$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("Test Data")
$ws3.Activate()
$excel.ActiveWindow.ScrollColumn = 11
$ws3.Range("B2:X3").Select() | Out-Null
Write-Host ("done: " + $excel.ActiveWindow.ScrollColumn)
